# Insert a narrow blank "spacer" column before each repeated report block
# (columns I, Q, Y, AG in the original layout). Inserting from right to
# left keeps the column indices for the remaining inserts stable.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$insertPositions = @(33, 25, 17, 9)
foreach ($p in $insertPositions) {
    $ws.Columns.Item($p).Insert()
    $ws.Columns.Item($p).ColumnWidth = 1.14
}

# Row 1 carries a bold header style on every cell, so the freshly
# inserted spacer columns pick up that style on their header row even
# though they hold no content. Fully clear those header cells so the
# spacer columns stay completely blank in row 1, same as every other row.
$headerGapCells = @("I1", "R1", "AA1", "AJ1")
foreach ($addr in $headerGapCells) {
    $ws.Range($addr).Clear()
}

